$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update password values for Noor.Uddin rows (5,6,7) from MHRA123456 to MHRA12345A
$ws.Range("B5").Value = "MHRA12345A"
$ws.Range("B6").Value = "MHRA12345A"
$ws.Range("B7").Value = "MHRA12345A"

# Update the selected cell on the sheet
$ws.Range("C17").Select()
